$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "239.03"

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "21.91"

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.452"

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.05651"

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "6.484"

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "3.352"

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "1.076"

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.7892"

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.1395"

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07396"

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.03208"

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.02973"

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.09255"

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.001664"

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.254"

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.04774"

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.0005740"

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.006237"

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.005117"

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.001049"

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.0001500"

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "3.892"

$ws.Range("E23").Value = "22LEOLEO"

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.1301"

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0004010"

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.04135"

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.006943"

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1043"

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.003010"

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.009905"

$ws.Range("E44").Value = "43LocalTradersLCTBestin24h"

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.00005435"

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00000000750"

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.6752"

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.03756"

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.00002100"

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.01010"
